# update project, add new function
#
# "Create Cylinder" gains a new "verificationDate" column (inserted right
# before "brand"), the TEST codes are bumped from ...001-003 to ...004-006,
# "1 van" is re-cased to "1 Van", and the sheet's selection moves onto the
# code column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create Cylinder")

# 1. Preserve the existing "brand" column (F) - values, shared-string refs
#    and styles - by copying it into the new column G before F is reused
#    for the new data.
$ws.Range("F1:F4").Copy($ws.Range("G1:G4")) | Out-Null

# 2. Turn the old "brand" header cell (F1) into the new "verificationDate"
#    header.
$ws.Range("F1").Value = "verificationDate"

# 3. Small data fixups already present on the sheet.
$ws.Range("D4").Value = "1 Van"
$ws.Range("A2").Value = "TEST1608004"
$ws.Range("A3").Value = "TEST1608005"
$ws.Range("A4").Value = "TEST1608006"

# 4. Fill the new verificationDate column with real dates, formatted with
#    the standard short-date number format (built-in numFmtId 14).
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Value = 44072
$ws.Range("F2").Copy($ws.Range("F3:F4")) | Out-Null
$ws.Range("F3").Value = 44075
$ws.Range("F4").Value = 44078

# 5. Approximate the visual column widths for the new/shifted columns -
#    verificationDate (F) takes on the width neighbouring column "weight"
#    (E) used, and brand (G) keeps the width the old "brand" column (F) had.
$ws.Columns.Item(6).ColumnWidth = 15.67
$ws.Columns.Item(7).ColumnWidth = 17.67

# 6. Match the new selection recorded for the sheet.
$ws.Range("A2:A4").Select() | Out-Null
